$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the leading '.' from the domain names in row 1 (B1:E1)
$ws.Range("B1").Value = "google.com"
$ws.Range("C1").Value = "nu.nl"
$ws.Range("D1").Value = "tweaker.net"
$ws.Range("E1").Value = "bbc.com"

# Update the active cell selection on the sheet (was D9, now E2)
$ws.Range("E2").Select()
